# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# (and one row on 演出) to reflect the latest scrape counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3333
$ws1.Range("F6").Value = 1094
$ws1.Range("F7").Value = 2226
$ws1.Range("F8").Value = 2138
$ws1.Range("F13").Value = 401
$ws1.Range("F15").Value = 45
$ws1.Range("F17").Value = 229
$ws1.Range("F19").Value = 644
$ws1.Range("F21").Value = 618
$ws1.Range("F22").Value = 12299
$ws1.Range("F23").Value = 12358
$ws1.Range("F29").Value = 383
$ws1.Range("F33").Value = 205
$ws1.Range("F34").Value = 600

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 38

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 3333
$ws4.Range("F7").Value = 1094
$ws4.Range("F8").Value = 2226
$ws4.Range("F9").Value = 2138
$ws4.Range("F14").Value = 401
$ws4.Range("F17").Value = 45
$ws4.Range("F21").Value = 229
$ws4.Range("F23").Value = 644
$ws4.Range("F25").Value = 618
$ws4.Range("F26").Value = 12299
$ws4.Range("F27").Value = 12358
$ws4.Range("F33").Value = 383
$ws4.Range("F39").Value = 205
$ws4.Range("F40").Value = 600
$ws4.Range("F41").Value = 38

$wb.Save()
